# Applies the cell-level text updates for the cryptos worksheet.
# All target cells hold plain text (inline strings in the original OOXML),
# so we force NumberFormat to Text ("@") before writing each value. That
# stops Excel's COM layer from auto-coercing number-looking strings (e.g.
# "0.998", "64.209.11") into floating point numbers, which would both change
# the cell type and mangle values like thousand-dot formatted prices.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '64.209.11'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  -1.49%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.106.63'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  -26.02%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.998'
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  +0.04%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '592.17'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  -1.25%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '157.35'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  +4.10%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.998'
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  -1.01%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.541'
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  -0.54%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '3.102.60'
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  -4.20%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.159'
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  -4.27%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '5.93'
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  -3.89%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.454'
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  -3.84%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.0000241'
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  -4.91%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '37.22'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  -4.65%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.120'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  -1.79%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '3.599.74'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  -2.88%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '7.24'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  -2.18%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '63.891.59'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  -1.35%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '3.089.49'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  -2.84%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '474.83'
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  -1.60%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '14.45'
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  -3.14%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.711'
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  -7.28%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '7.59'
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  -3.97%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.45'
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  -1.98%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '12.97'
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  -5.28%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '80.95'
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  -3.04%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '10.41'
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  +3.64%  '
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  -0.11%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '7.48'
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  +0.40%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '2.68'
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  -3.71%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.998'
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  +0.16%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '2.18'
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  -4.89%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.113'
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  -6.49%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '27.41'
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  -5.35%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.0₃0845'
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  -3.40%  '
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  -2.90%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '6.05'
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  -4.85%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '3.30'
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  -2.44%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.24'
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  -5.60%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '50.95'
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  -2.61%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '9.18'
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  -3.24%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '438.89'
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  -8.51%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.291'
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  -3.71%  '
$ws.Range('B44').NumberFormat = '@'
$ws.Range('B44').Value = 'Arweave'
$ws.Range('C44').NumberFormat = '@'
$ws.Range('C44').Value = 'https://coinranking.com/coin/7XWg41D1+arweave-ar'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '40.38'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  +1.96%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.0363'
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  -4.84%  '
$ws.Range('B46').NumberFormat = '@'
$ws.Range('B46').Value = 'Kaspa'
$ws.Range('C46').NumberFormat = '@'
$ws.Range('C46').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.112'
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  -0.67%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '2.824.30'
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  -4.21%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '130.08'
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  -1.68%  '
$ws.Range('B49').NumberFormat = '@'
$ws.Range('B49').Value = 'InjectiveProtocol'
$ws.Range('C49').NumberFormat = '@'
$ws.Range('C49').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '25.38'
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  -0.34%  '
$ws.Range('B50').NumberFormat = '@'
$ws.Range('B50').Value = 'USDe'
$ws.Range('C50').NumberFormat = '@'
$ws.Range('C50').Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.999'
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  +0.00%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '2.24'
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  -3.80%  '
